# Add "Devin/Paul" as a name to the task-assignment list on the
# "Improvements" sheet (row 6: "Redesign the class UML."), and make the
# "Improvements" tab the active/selected sheet (it was "Bug Fixes" before).

$wb = $excel.ActiveWorkbook

# --- Improvements sheet: add the new contributor to row 6 ---------------
$ws = $wb.Worksheets.Item("Improvements")

# Column C gets the newly-added name, column D reuses the existing
# "MarchMadnessGUI" file reference (same text used elsewhere in the sheet).
$ws.Range("C6").Value = "Devin/Paul"
$ws.Range("D6").Value = "MarchMadnessGUI"

# --- Switch the active sheet/selection -----------------------------------
# Previously "Bug Fixes" was the selected/active tab with cell C11
# selected; now "Improvements" becomes the active tab with D6 selected.
$ws.Activate()
$ws.Range("D6").Select()
